$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last updated" timestamp title in A1
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 15:20"

# 2) Serbia moves up in the ranking (new case counts put it ahead of Eslovenia).
#    Insert a fresh row just above Eslovenia's current row and populate it with
#    Serbia's updated figures, then remove Serbia's old row further down the list.
$eslovenia = $ws.Range("A:A").Find("Eslovenia")
$esloveniaRow = $eslovenia.Row

$ws.Rows($esloveniaRow).Insert()

$newRow = $esloveniaRow
$ws.Range("A" + $newRow).Value = "Serbia"
$ws.Range("B" + $newRow).Value = 741
$ws.Range("C" + $newRow).Value = 82
$ws.Range("D" + $newRow).Value = 42
$ws.Range("E" + $newRow).Value = 686
$ws.Range("F" + $newRow).Value = 25
$ws.Range("G" + $newRow).Value = 3
$ws.Range("H" + $newRow).Value = 13

$searchBelow = $ws.Range("A" + ($newRow + 1) + ":A1000")
$oldSerbia = $searchBelow.Find("Serbia")
$ws.Rows($oldSerbia.Row).Delete()

# 3) Refresh case numbers for a handful of countries.
function Set-CountryRow($countryName, $totalCases, $newCases, $activeCases, $recovered, $critical, $deathsToday, $deaths) {
    $found = $ws.Range("A:A").Find($countryName)
    $r = $found.Row
    $ws.Range("B" + $r).Value = $totalCases
    $ws.Range("C" + $r).Value = $newCases
    $ws.Range("D" + $r).Value = $activeCases
    $ws.Range("E" + $r).Value = $recovered
    $ws.Range("F" + $r).Value = $critical
    $ws.Range("G" + $r).Value = $deathsToday
    $ws.Range("H" + $r).Value = $deaths
}

Set-CountryRow "Reino Unido" 19522 2433 135 18159 163 209 1228
Set-CountryRow "Suiza" 14593 517 1595 12708 301 26 290
Set-CountryRow "Austria" 8536 265 479 7971 187 18 86
Set-CountryRow "Arabia Saudita" 1299 96 66 1225 12 4 8
Set-CountryRow "Finlandia" 1239 72 10 1218 32 2 11
Set-CountryRow "Republica de Macedonia" 259 18 3 250 1 2 6
